$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 12545
$ws.Range("I6").Value = 16693.334
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 50080.00199999999
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -49968.00199999999
$ws.Range("N6").Value = -524
# Row 9
$ws.Range("H9").Value = 138.38461
$ws.Range("I9").Value = 138.63637
$ws.Range("J9").Value = 137
$ws.Range("K9").Value = 138.63637
$ws.Range("L9").Value = 137
$ws.Range("M9").Value = 30.36363
$ws.Range("N9").Value = -475
# Row 12
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -640
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 93
$ws.Range("H93").Value = 20601
$ws.Range("J93").Value = 20601
$ws.Range("L93").Value = 20601
$ws.Range("N93").Value = -25593
# Row 113
$ws.Range("H113").Value = 20003390
$ws.Range("I113").Value = 50002476
$ws.Range("K113").Value = 50002476
$ws.Range("M113").Value = -49999222
# Row 132
$ws.Range("H132").Value = 8339787
$ws.Range("I132").Value = 9806632
$ws.Range("J132").Value = 27666.666
$ws.Range("K132").Value = 29419896
$ws.Range("L132").Value = 82999.99800000001
$ws.Range("M132").Value = -29417366
$ws.Range("N132").Value = -88059.99800000001
# Row 137
$ws.Range("H137").Value = 1043.9
$ws.Range("I137").Value = 1022.56
$ws.Range("K137").Value = 3067.68
$ws.Range("M137").Value = -517.6799999999998
# Row 141
$ws.Range("H141").Value = 895
$ws.Range("I141").Value = 895
$ws.Range("K141").Value = 2685
$ws.Range("M141").Value = 2495

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 62
$ws.Range("H62").Value = 69812
$ws.Range("J62").Value = 69812
$ws.Range("L62").Value = 69812
$ws.Range("N62").Value = -71060
# Row 65
$ws.Range("H65").Value = 69812
$ws.Range("J65").Value = 69812
$ws.Range("L65").Value = 209436
$ws.Range("N65").Value = -215676
# Row 97
$ws.Range("H97").Value = 487.44446
$ws.Range("J97").Value = 899
$ws.Range("L97").Value = 899
$ws.Range("N97").Value = -1891

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 17857978
$ws.Range("I94").Value = 19231516
$ws.Range("K94").Value = 19231516
$ws.Range("M94").Value = -19231065

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1808.238
$ws.Range("I31").Value = 865.86365
$ws.Range("K31").Value = 865.86365
$ws.Range("M31").Value = -570.86365
# Row 34
$ws.Range("H34").Value = 1808.238
$ws.Range("I34").Value = 865.86365
$ws.Range("K34").Value = 865.86365
$ws.Range("M34").Value = -663.86365
# Row 35
$ws.Range("H35").Value = 152.5
$ws.Range("I35").Value = 152.5
$ws.Range("K35").Value = 152.5
$ws.Range("M35").Value = 141.5
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()
# Row 62
$ws.Range("H62").Value = 7695235
$ws.Range("I62").Value = 3044.4
$ws.Range("J62").Value = 200000000
$ws.Range("K62").Value = 3044.4
$ws.Range("L62").Value = 200000000
$ws.Range("M62").Value = -2420.4
$ws.Range("N62").Value = -200001248
# Row 65
$ws.Range("H65").Value = 7695235
$ws.Range("I65").Value = 3044.4
$ws.Range("J65").Value = 200000000
$ws.Range("K65").Value = 15222
$ws.Range("L65").Value = 1000000000
$ws.Range("M65").Value = -12102
$ws.Range("N65").Value = -1000006240
# Row 122
$ws.Range("H122").Value = 687.5
$ws.Range("I122").Value = 687.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2062.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 387.5
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 7687.222
$ws.Range("I132").Value = 11137.2
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 33411.60000000001
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -30881.60000000001
$ws.Range("N132").Value = -15184.25

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 3778.6428
$ws.Range("I63").Value = 871.5714
$ws.Range("J63").Value = 6685.7144
$ws.Range("K63").Value = 2614.7142
$ws.Range("L63").Value = 20057.1432
$ws.Range("M63").Value = -1865.7142
$ws.Range("N63").Value = -21555.1432
# Row 66
$ws.Range("H66").Value = 3778.6428
$ws.Range("I66").Value = 871.5714
$ws.Range("J66").Value = 6685.7144
$ws.Range("K66").Value = 7844.1426
$ws.Range("L66").Value = 60171.4296
$ws.Range("M66").Value = -4100.1426
$ws.Range("N66").Value = -67659.4296
# Row 131
$ws.Range("H131").Value = 41668524
$ws.Range("I131").Value = 250000350
$ws.Range("J131").Value = 2158.9
$ws.Range("K131").Value = 750001050
$ws.Range("L131").Value = 6476.700000000001
$ws.Range("M131").Value = -749996010
$ws.Range("N131").Value = -16556.7

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 109
$ws.Range("H109").Value = 29198.2
$ws.Range("I109").Value = 28997
$ws.Range("J109").Value = 29500
$ws.Range("K109").Value = 28997
$ws.Range("L109").Value = 29500
$ws.Range("M109").Value = -27957
$ws.Range("N109").Value = -31580
# Row 113
$ws.Range("H113").Value = 2286.6365
$ws.Range("I113").Value = 1365.375
$ws.Range("J113").Value = 2813.0715
$ws.Range("K113").Value = 1365.375
$ws.Range("L113").Value = 2813.0715
$ws.Range("M113").Value = 804.625
$ws.Range("N113").Value = -7153.0715

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1149.3334
$ws.Range("I16").Value = 1055.0834
$ws.Range("J16").Value = 1526.3334
$ws.Range("K16").Value = 1055.0834
$ws.Range("L16").Value = 1526.3334
$ws.Range("M16").Value = -885.0834
$ws.Range("N16").Value = -1866.3334
# Row 61
$ws.Range("H61").Value = 1214.6364
$ws.Range("I61").Value = 1239.5555
$ws.Range("J61").Value = 1102.5
$ws.Range("K61").Value = 1239.5555
$ws.Range("L61").Value = 1102.5
$ws.Range("M61").Value = -1037.5555
$ws.Range("N61").Value = -1506.5
# Row 68
$ws.Range("H68").Value = 1732.125
$ws.Range("I68").Value = 1654.9231
$ws.Range("K68").Value = 1654.9231
$ws.Range("M68").Value = -905.9231
# Row 71
$ws.Range("H71").Value = 1732.125
$ws.Range("I71").Value = 1654.9231
$ws.Range("K71").Value = 8274.6155
$ws.Range("M71").Value = -4530.6155
# Row 113
$ws.Range("H113").Value = 1214.6364
$ws.Range("I113").Value = 1239.5555
$ws.Range("J113").Value = 1102.5
$ws.Range("K113").Value = 1239.5555
$ws.Range("L113").Value = 1102.5
$ws.Range("M113").Value = 930.4445000000001
$ws.Range("N113").Value = -5442.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2523.2432
$ws.Range("I132").Value = 1888.8387
$ws.Range("K132").Value = 5666.5161
$ws.Range("M132").Value = -3136.5161
